$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 46081 (2026-02-28)
# to 46082 (2026-03-01) for every data row (rows 2 through 161).
$lastRow = 161
$ws.Range("C2:C$lastRow").Value = 46082
